$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update (A1) ---
$ws.Range("A1").Value() = "Datos actualizados a 20 de Septiembre de 2020 a las 12:16"

# --- Swap country rows: Timor Oriental (row 204) <-> Santa Lucia (row 205) ---
$ws.Range("A204").Value() = "Santa Lucia"
$ws.Range("A205").Value() = "Timor Oriental"

# --- Swap country rows: Islas Malvinas (row 214) <-> Montserrat (row 215) ---
$ws.Range("A214").Value() = "Montserrat"
$ws.Range("A215").Value() = "Islas Malvinas"

# Data for these two rows also swaps along with the country identity
# (row 214 was Islas Malvinas w/ D=13,H=0 ; row 215 was Montserrat w/ D=12,H=1)
$ws.Range("D214").Value() = 12
$ws.Range("H214").Value() = 1
$ws.Range("D215").Value() = 13
$ws.Range("H215").Value() = 0

# --- Updated covid stats (provincias/paises refresh) ---

# Row 5: India
$ws.Range("B5").Value() = 5405252
$ws.Range("C5").Value() = 7022
$ws.Range("E5").Value() = 1015413
$ws.Range("G5").Value() = 22
$ws.Range("H5").Value() = 86796

# Row 18: Banglades
$ws.Range("B18").Value() = 348918
$ws.Range("C18").Value() = 1544
$ws.Range("D18").Value() = 256565
$ws.Range("E18").Value() = 87414
$ws.Range("G18").Value() = 26
$ws.Range("H18").Value() = 4939

# Row 25: Alemania
$ws.Range("B25").Value() = 272311
$ws.Range("C25").Value() = 3
$ws.Range("E25").Value() = 19345

# Row 46: Emiratos Arabes Unidos
$ws.Range("B46").Value() = 84916
$ws.Range("C46").Value() = 674
$ws.Range("D46").Value() = 74273
$ws.Range("E46").Value() = 10239

# Row 76: El Salvador
$ws.Range("B76").Value() = 27553
$ws.Range("C76").Value() = 125
$ws.Range("D76").Value() = 21561
$ws.Range("E76").Value() = 5181
$ws.Range("G76").Value() = 3
$ws.Range("H76").Value() = 811

# Row 103: Finlandia
$ws.Range("B103").Value() = 8980
$ws.Range("C103").Value() = 58
$ws.Range("E103").Value() = 941
